$wb = $excel.ActiveWorkbook

# ---- C48: fill rows 122-200 (point group + energy) ----
$ws48 = $wb.Worksheets.Item("C48")
$data48 = @(
    @(122, "C1", 25.257),
    @(123, "C2", 26.92),
    @(124, "C1", 26.905),
    @(125, "C1", 24.697),
    @(126, "C2", 25.085),
    @(127, "C1", 24.594),
    @(128, "C1", 24.052),
    @(129, "C1", 24.637),
    @(130, "C1", 24.642),
    @(131, "C1", 24.793),
    @(132, "C1", 24.639),
    @(133, "C1", 25.68),
    @(134, "C1", 24.258),
    @(135, "C2", 25.66),
    @(136, "C1", 23.629),
    @(137, "Cs", 25.144),
    @(138, "C2", 23.48),
    @(139, "C2v", 23.526),
    @(140, "C1", 23.171),
    @(141, "C1", 23.351),
    @(142, "Cs", 24.086),
    @(143, "Cs", 24.077),
    @(144, "C1", 25.055),
    @(145, "D2", 26.824),
    @(146, "C1", 23.52),
    @(147, "C1", 24.162),
    @(148, "C1", 23.756),
    @(149, "Cs", 23.347),
    @(150, "C1", 22.944),
    @(151, "C1", 23.231),
    @(152, "C1", 23.839),
    @(153, "C2", 24.375),
    @(154, "C2", 24.445),
    @(155, "C1", 24.044),
    @(156, "C1", 23.56),
    @(157, "C1", 23.483),
    @(158, "C2", 25.159),
    @(159, "C2", 23.965),
    @(160, "C1", 23.999),
    @(161, "C2", 23.765),
    @(162, "C2", 23.833),
    @(163, "C1", 23.543),
    @(164, "C2", 23.132),
    @(165, "C1", 23.948),
    @(166, "C1", 22.927),
    @(167, "C2", 23.565),
    @(168, "C1", 23.288),
    @(169, "Cs", 23.487),
    @(170, "D2", 23.276),
    @(171, "C2", 22.624),
    @(172, "C2", 23.101),
    @(173, "C1", 24.645),
    @(174, "C1", 23.688),
    @(175, "C1", 24.792),
    @(176, "C2", 24.138),
    @(177, "C1", 23.717),
    @(178, "C2", 24),
    @(179, "C1", 23.628),
    @(180, "C1", 24.1),
    @(181, "C2", 24.037),
    @(182, "C1", 25.528),
    @(183, "C1", 25.046),
    @(184, "C2", 23.468),
    @(185, "Cs", 24.292),
    @(186, "C2", 23.381),
    @(187, "D6d", 26.953),
    @(188, "Cs", 23.245),
    @(189, "D3", 22.993),
    @(190, "D6d", 26.643),
    @(191, "C2", 24.979),
    @(192, "C2", 24.796),
    @(193, "C2", 23.096),
    @(194, "C1", 23.436),
    @(195, "C2", 23.104),
    @(196, "C1", 22.944),
    @(197, "C1", 24.718),
    @(198, "Cs", 22.946),
    @(199, "D2", 23.815),
    @(200, "C2", 22.376)
)
foreach ($row in $data48) {
    $ws48.Cells.Item($row[0], 2).Value = $row[1]
    $ws48.Cells.Item($row[0], 3).Value = $row[2]
}

# ---- C50: fill rows 2-9 (point group + energy) ----
$ws50 = $wb.Worksheets.Item("C50")
$data50 = @(
    @(2, "D5h", 34.883),
    @(3, "C2", 31.919),
    @(4, "D3h", 35.747),
    @(5, "Cs", 28.936),
    @(6, "Cs", 29.08),
    @(7, "C2", 30.333),
    @(8, "C1", 29.243),
    @(9, "Cs", 28.39)
)
foreach ($row in $data50) {
    $ws50.Cells.Item($row[0], 2).Value = $row[1]
    $ws50.Cells.Item($row[0], 3).Value = $row[2]
}

# ---- Row height quirk: row 122 on C48 ----
$ws48.Rows.Item(122).RowHeight = 12.1

# ---- View/selection state ----
$ws48.Range("B201").Select()
$ws50.Activate()
$ws50.Range("B10").Select()